# Daily attendance processing - 2026-01-06 23:02:34
# Reorders the comma-separated "Recorded By" names in column G so that any
# "System" / "system" token(s) come first, followed by the remaining
# tokens (e.g. email addresses) in their original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($raw -eq $null) {
        continue
    }

    $text = [string]$raw
    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $parts = $text.Split(",")

    $systemParts = @()
    $otherParts = @()

    foreach ($p in $parts) {
        $trimmed = $p.Trim()
        if ($trimmed.ToLower() -eq "system") {
            $systemParts = $systemParts + $trimmed
        } else {
            $otherParts = $otherParts + $trimmed
        }
    }

    if ($systemParts.Length -eq 0) {
        continue
    }

    $newParts = $systemParts + $otherParts
    $newText = $newParts -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
